$d = $word.ActiveDocument

# 1) Merge the "{% if not loop.last %} " paragraph with the literal
#    "<w:p>...<w:br.../>...</w:p>" paragraph by deleting the paragraph
#    mark between them and the trailing space before it.
$d.Content.Find.Execute("%} " + [char]13 + "<", $false, $false, $false, $false, $false, $true, 1, $false, "%}<", 2) | Out-Null

# 2) Merge that paragraph with the "{% endif %} " paragraph the same way.
$d.Content.Find.Execute("</w:p> " + [char]13 + "{", $false, $false, $false, $false, $false, $true, 1, $false, "</w:p>{", 2) | Out-Null

# 3) Remove the old _GoBack bookmark (currently sitting at the start of the
#    now-merged "{% if not loop.last %}..." paragraph).
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# 4) Re-create the _GoBack bookmark immediately after "{% endfor %}" (at the
#    end of the last paragraph, before its paragraph mark). A temporary
#    marker character is used because adding a bookmark directly from a
#    zero-width Range placed exactly before a paragraph mark is unreliable;
#    routing the placement through Find (non-collapsed match) and then
#    deleting the marker lands the bookmark correctly.
$marker = [char]1
$d.Content.Find.Execute("{% endfor %}", $false, $false, $false, $false, $false, $true, 1, $false, "{% endfor %}" + $marker, 2) | Out-Null

$markerRange = $d.Content
$markerRange.Find.Execute($marker, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("_GoBack", $markerRange)

$d.Content.Find.Execute($marker, $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
